$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data-driven cell updates derived from the commit diff: each tuple is (CellRef, NewValue)
$updates = @(
    @('D2', '30.553.28'),
    @('E2', '  +0.50%  '),
    @('D3', '1.916.83'),
    @('E3', '  -0.15%  '),
    @('E4', '  +0.19%  '),
    @('D5', '243.90'),
    @('E5', '  +1.22%  '),
    @('E6', '  +0.11%  '),
    @('D7', '0.4925'),
    @('E7', '  +4.98%  '),
    @('D8', '0.2901'),
    @('E8', '  +1.47%  '),
    @('D9', '0.06707'),
    @('E9', '  -3.40%  '),
    @('B10', 'Solana'),
    @('C10', 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'),
    @('D10', '18.89'),
    @('E10', '  +3.08%  '),
    @('B11', 'Litecoin'),
    @('C11', 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'),
    @('D11', '106.25'),
    @('E11', '  -0.79%  '),
    @('E12', '  +0.23%  '),
    @('D13', '0.07643'),
    @('E13', '  -0.10%  '),
    @('D14', '5.251'),
    @('E14', '  +1.24%  '),
    @('D15', '0.6655'),
    @('E15', '  +1.13%  '),
    @('D16', '272.56'),
    @('E16', '  -7.24%  '),
    @('D17', '30.534.68'),
    @('E17', '  +0.41%  '),
    @('E18', '  +0.06%  '),
    @('D19', '0.000007522'),
    @('E19', '  -1.66%  '),
    @('D20', '2.169.67'),
    @('E20', '  +0.27%  '),
    @('E21', '  -1.68%  '),
    @('D22', '5.481'),
    @('E22', '  +4.95%  '),
    @('E23', '  -0.21%  '),
    @('D24', '6.392'),
    @('E24', '  +2.95%  '),
    @('D25', '9.398'),
    @('E25', '  +1.45%  '),
    @('D26', '164.49'),
    @('E26', '  -2.06%  '),
    @('D27', '19.81'),
    @('E27', '  -7.38%  '),
    @('D28', '2.098'),
    @('E28', '  +2.54%  '),
    @('D29', '0.1053'),
    @('E29', '  -2.51%  '),
    @('D30', '1.408'),
    @('E30', '  +3.14%  '),
    @('D31', '4.133'),
    @('E31', '  -0.45%  '),
    @('D32', '4.024'),
    @('E32', '  +1.31%  '),
    @('D33', '0.04989'),
    @('E33', '  -1.29%  '),
    @('D34', '0.7250'),
    @('E34', '  -2.66%  '),
    @('D35', '1.132'),
    @('E35', '  -1.17%  '),
    @('D36', '1.000'),
    @('E36', '  +0.05%  '),
    @('E37', '  -0.19%  '),
    @('D38', '0.02030'),
    @('E38', '  +0.53%  '),
    @('E39', '  -0.27%  '),
    @('D40', '110.95'),
    @('E40', '  +2.30%  '),
    @('E41', '  -1.82%  '),
    @('D42', '0.4398'),
    @('E42', '  +4.34%  '),
    @('D43', '0.8671'),
    @('E43', '  -0.53%  '),
    @('D44', '5.863'),
    @('E44', '  +0.25%  '),
    @('D45', '1.001'),
    @('E45', '  +0.11%  '),
    @('D46', '67.82'),
    @('E46', '  +0.34%  '),
    @('D47', '7.237'),
    @('E47', '  +0.69%  '),
    @('D48', '9.257'),
    @('E48', '  +0.81%  '),
    @('D49', '0.1245'),
    @('E49', '  +2.91%  '),
    @('D50', '47.79'),
    @('E50', '  -11.54%  '),
    @('D51', '34.62'),
    @('E51', '  -0.14%  ')
)

foreach ($u in $updates) {
    $cellRef = $u[0]
    $newVal = $u[1]
    $rng = $ws.Range($cellRef)
    # Force text format so numeric-looking strings (e.g. "30.553.28", "1.000")
    # are preserved verbatim instead of being coerced into Excel numbers/dates.
    $rng.NumberFormat = "@"
    $rng.Value = $newVal
}
